$d = $word.ActiveDocument

# --- Change 1 -------------------------------------------------------------
# "Add 1 for each strenuous action performed. Then, adjust for encumbrance
#  and worn armour (see page XXX)." ->
# "Add 1 for each strenuous action performed (note: This can be modified
#  if you're wearing medium or heavier armour)."
$apos = [char]0x2019
$newSentence = "Add 1 for each strenuous action performed (note: This can be modified if you" + $apos + "re wearing medium or heavier armour)."

$rng1 = $d.Content
$found1 = $rng1.Find.Execute(
    "Add 1 for each strenuous action performed. Then, adjust for encumbrance and worn armour (see page XXX).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    $newSentence, 2)

if (-not $found1) {
    throw "Could not find the fatigue-test sentence to replace."
}

# --- Change 2 ---------------------------------------------------------------
# The red "<Max 2>" marker becomes "<Max 3>" (still red).
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Max 2", $true, $false, $false, $false, $false, $true, 1, $false, "Max 3", 2)

if (-not $found2) {
    throw "Could not find the '<Max 2>' marker to replace."
}

Write-Host "Change 1 applied: $found1"
Write-Host "Change 2 applied: $found2"
